# feat: add transaction 00-0267
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell AG1 mirrors the existing "Transacción" header (column S1)
$ws.Range("AG1").Value = "Transacción"

# New data cell AG2 carries the new transaction value for this row
$ws.Range("AG2").Value = "00-0267"

# Selection moves back to the top-left of the sheet
$ws.Range("A2").Select()
